$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.331.90'
$ws.Range('E2').Value = '  -1.39%  '
$ws.Range('D3').Value = '2.367.41'
$ws.Range('E3').Value = '  +4.59%  '
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.60'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.648'
$ws.Range('E6').Value = '  -1.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '72.05'
$ws.Range('E7').Value = '  +12.96%  '
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.479'
$ws.Range('E9').Value = '  +6.29%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0972'
$ws.Range('E10').Value = '  -0.61%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '56.87'
$ws.Range('E11').Value = '  -2.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '27.26'
$ws.Range('E12').Value = '  +2.52%  '
$ws.Range('D13').Value = '2.729.14'
$ws.Range('E13').Value = '  +4.97%  '
$ws.Range('E14').Value = '  +0.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.95'
$ws.Range('E15').Value = '  +2.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.27'
$ws.Range('E16').Value = '  +2.32%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.859'
$ws.Range('E17').Value = '  +2.10%  '
$ws.Range('D18').Value = '2.381.04'
$ws.Range('E18').Value = '  +4.96%  '
$ws.Range('D19').Value = '43.367.53'
$ws.Range('E19').Value = '  -1.08%  '
$ws.Range('D20').Value = '0.0₃0992'
$ws.Range('E20').Value = '  +0.98%  '
$ws.Range('B21').Value = 'Litecoin'
$ws.Range('C21').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '74.50'
$ws.Range('E21').Value = '  +0.92%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.33'
$ws.Range('E22').Value = '  +2.51%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '249.29'
$ws.Range('E23').Value = '  -0.24%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  +0.00%  '
$ws.Range('E25').Value = '  +3.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.45'
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('E27').Value = '  +1.09%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.97'
$ws.Range('E28').Value = '  +0.76%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '22.43'
$ws.Range('E29').Value = '  +2.27%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '173.85'
$ws.Range('E30').Value = '  -0.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.51'
$ws.Range('E31').Value = '  +5.49%  '
$ws.Range('E32').Value = '  -5.66%  '
$ws.Range('E33').Value = '  -0.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.97'
$ws.Range('E34').Value = '  -0.24%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0692'
$ws.Range('E35').Value = '  +0.97%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.04'
$ws.Range('E36').Value = '  +1.48%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.59'
$ws.Range('E37').Value = '  +2.49%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.70'
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('B39').Value = 'LidoDAOToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.44'
$ws.Range('E39').Value = '  +6.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0254'
$ws.Range('E40').Value = '  -0.40%  '
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.88'
$ws.Range('E42').Value = '  +0.88%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '18.42'
$ws.Range('E43').Value = '  +6.00%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.18'
$ws.Range('E44').Value = '  +8.72%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '99.85'
$ws.Range('E45').Value = '  +1.10%  '
$ws.Range('B46').Value = 'FTXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.51'
$ws.Range('E46').Value = '  -1.25%  '
$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.21'
$ws.Range('E47').Value = '  +0.97%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0953'
$ws.Range('E48').Value = '  +0.10%  '
$ws.Range('D49').Value = '1.442.84'
$ws.Range('E49').Value = '  -0.85%  '
$ws.Range('D50').Value = '2.600.21'
$ws.Range('E50').Value = '  +5.07%  '
$ws.Range('E51').Value = '  -2.61%  '
